# Scénario de test v1 - Première version des scénarii de tests.
# Renseigne la colonne "RESULTAT OBSERVE" (F) pour les cas de test 6.1/6.2 (lignes 13-15)
# et met en évidence la ligne 14 (cas 6.2, crash applicatif) avec la mise en forme
# "erreur" déjà utilisée ailleurs dans le classeur.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Copie la mise en forme "erreur" (fond rouge) de la ligne 4 vers la ligne 14
# (colonnes UTILISATEUR, LIBELLE TEST, RESULTAT ATTENDU, RESULTAT OBSERVE)
# sans toucher au contenu des cellules.
$ws.Range("C4:F4").Copy()
$ws.Range("C14:F14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Renseigne les résultats observés.
$ws.Range("F13").Value = "OK"
$ws.Range("F14").Value = "A revérifier - crash appli - issue ouverte"
$ws.Range("F15").Value = "OK"

# Met à jour la sélection active de la feuille.
$ws.Range("E14").Select()
